$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 443, shifting the existing rows 443:458 down
# to 444:459 (dimension grows from A1:R458 to A1:R459).
$ws.Rows(443).Insert()

# Populate the newly inserted row 443 with the latest weekly price record.
$ws.Range("A443").Value = 4
$ws.Range("B443").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C443").Value = "Los Lagos"
$ws.Range("D443").Value = 45239
$ws.Range("E443").Value = 10
$ws.Range("F443").Value = 100112028
$ws.Range("G443").Value = "Sandia"
$ws.Range("H443").Value = "Sin especificar"
$ws.Range("I443").Value = "Primera"
$ws.Range("J443").Value = 500
$ws.Range("K443").Value = 1300
$ws.Range("L443").Value = 1300
$ws.Range("M443").Value = 1300
$ws.Range("N443").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O443").Value = "Perú"
$ws.Range("P443").Value = 1300
$ws.Range("Q443").Value = 1
$ws.Range("R443").Value = "Hortaliza"
